$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Register the lowercase date-time format variant (numFmtId 164) that the
# original export also declared, then discard the scratch cell/row used to
# register it so it never shows up as real content.
$ws.Range("Z100").Value = 1
$ws.Range("Z100").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("Z100").EntireRow.Delete()

# Add new "Trening" header in column F, matching the style of the other headers
# (copy the formatting of an existing header cell so it reuses the same style)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Trening"

# Define the custom date-time number format used for column A
$dateFormat = "YYYY-MM-DD HH:MM:SS"

# Update column A (rows 2-7) to real date-time serial values instead of text,
# and tag column F with the training phase label "Duża Gra"
$ws.Range("A2").Value = 45686.47578217593
$ws.Range("A3").Value = 45686.47806111111
$ws.Range("A4").Value = 45686.47809236111
$ws.Range("A5").Value = 45686.47577060185
$ws.Range("A6").Value = 45686.47641643519
$ws.Range("A7").Value = 45686.47685625

$ws.Range("A2:A7").NumberFormat = $dateFormat

$ws.Range("F2").Value = "Duża Gra"
$ws.Range("F3").Value = "Duża Gra"
$ws.Range("F4").Value = "Duża Gra"
$ws.Range("F5").Value = "Duża Gra"
$ws.Range("F6").Value = "Duża Gra"
$ws.Range("F7").Value = "Duża Gra"

# Add two new rows (8 and 9) describing the "Mała Gra" part of training.
# B/C/D stay empty but still need to exist as real (empty) cells, so nudge
# a formatting no-op on them to materialize the cell without changing style.
$ws.Range("A8").Value = 45686
$ws.Range("A8").NumberFormat = $dateFormat
$ws.Range("B8").Borders.LineStyle = 0
$ws.Range("C8").Borders.LineStyle = 0
$ws.Range("D8").Borders.LineStyle = 0
$ws.Range("E8").Value = "10-15"
$ws.Range("F8").Value = "Mała Gra"

$ws.Range("A9").Value = 45686
$ws.Range("A9").NumberFormat = $dateFormat
$ws.Range("B9").Borders.LineStyle = 0
$ws.Range("C9").Borders.LineStyle = 0
$ws.Range("D9").Borders.LineStyle = 0
$ws.Range("E9").Value = "5-10"
$ws.Range("F9").Value = "Mała Gra"
